# Auto-generated edit script: update sensor data rows and append new rows
# Shifts existing data up by one row (dropping the old row2 values) and
# appends 10 new rows (22-31) of walkingToRunning sensor data, extending
# the sheet from A1:H21 to A1:H31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 0.0
$ws.Range("B2").Value = "walkingToRunning"
$ws.Range("C2").Value = 10.44645118713379
$ws.Range("D2").Value = -1.072464942932129
$ws.Range("E2").Value = -7.741544246673584
$ws.Range("F2").Value = -0.4181592197061701
$ws.Range("G2").Value = 4.999134939407615
$ws.Range("H2").Value = 4.467728652686705

# Row 3
$ws.Range("A3").Value = 100.0
$ws.Range("B3").Value = "walkingToRunning"
$ws.Range("C3").Value = -12.54131031036377
$ws.Range("D3").Value = -6.995788097381592
$ws.Range("E3").Value = 4.960752010345459
$ws.Range("F3").Value = -0.2939812507584775
$ws.Range("G3").Value = 1.309815360006858
$ws.Range("H3").Value = 2.089699776373152

# Row 4
$ws.Range("A4").Value = 200.0
$ws.Range("B4").Value = "walkingToRunning"
$ws.Range("C4").Value = 1.563988208770752
$ws.Range("D4").Value = -0.6327500343322754
$ws.Range("E4").Value = 7.09266471862793
$ws.Range("F4").Value = 0.6610930461749918
$ws.Range("G4").Value = -2.38056925078418
$ws.Range("H4").Value = -1.609635165918661

# Row 5
$ws.Range("A5").Value = 300.0
$ws.Range("B5").Value = "walkingToRunning"
$ws.Range("C5").Value = 1.386142730712891
$ws.Range("D5").Value = -3.08348274230957
$ws.Range("E5").Value = 6.775681018829346
$ws.Range("F5").Value = 1.810891748588788
$ws.Range("G5").Value = -1.378563615763306
$ws.Range("H5").Value = -2.978032114349805

# Row 6
$ws.Range("A6").Value = 400.0
$ws.Range("B6").Value = "walkingToRunning"
$ws.Range("C6").Value = 1.111974596977234
$ws.Range("D6").Value = -5.867249965667725
$ws.Range("E6").Value = 13.13678932189941
$ws.Range("F6").Value = 2.31189947707626
$ws.Range("G6").Value = -2.220112597830942
$ws.Range("H6").Value = -3.509698769756517

# Row 7
$ws.Range("A7").Value = 500.0
$ws.Range("B7").Value = "walkingToRunning"
$ws.Range("C7").Value = -4.935175895690918
$ws.Range("D7").Value = 8.790708541870117
$ws.Range("E7").Value = 6.229645729064941
$ws.Range("F7").Value = -4.20486067611483
$ws.Range("G7").Value = -2.35995082988918
$ws.Range("H7").Value = 3.951543032566025

# Row 8
$ws.Range("A8").Value = 600.0
$ws.Range("B8").Value = "walkingToRunning"
$ws.Range("C8").Value = 16.58664321899414
$ws.Range("D8").Value = -65.17784118652344
$ws.Range("E8").Value = 6.808780193328857
$ws.Range("F8").Value = -6.209334919385752
$ws.Range("G8").Value = -1.759963089060472
$ws.Range("H8").Value = 10.33535447298917

# Row 9
$ws.Range("A9").Value = 700.0
$ws.Range("B9").Value = "walkingToRunning"
$ws.Range("C9").Value = -8.056964874267578
$ws.Range("D9").Value = 23.69245338439941
$ws.Range("E9").Value = -1.875571012496948
$ws.Range("F9").Value = 1.966778721764825
$ws.Range("G9").Value = 7.234104103016977
$ws.Range("H9").Value = 2.021741985160567

# Row 10
$ws.Range("A10").Value = 800.0
$ws.Range("B10").Value = "walkingToRunning"
$ws.Range("C10").Value = -4.308743476867676
$ws.Range("D10").Value = 2.427361011505127
$ws.Range("E10").Value = 6.903748512268066
$ws.Range("F10").Value = 2.757100457342962
$ws.Range("G10").Value = -3.825466089159538
$ws.Range("H10").Value = -3.918381726630348

# Row 11
$ws.Range("A11").Value = 900.0
$ws.Range("B11").Value = "walkingToRunning"
$ws.Range("C11").Value = -2.216989278793335
$ws.Range("D11").Value = -25.4277229309082
$ws.Range("E11").Value = 2.421021461486816
$ws.Range("F11").Value = 6.505795770716528
$ws.Range("G11").Value = -1.982219981255937
$ws.Range("H11").Value = -9.447731205236206

# Row 12
$ws.Range("A12").Value = 1000.0
$ws.Range("B12").Value = "walkingToRunning"
$ws.Range("C12").Value = -0.3425538539886474
$ws.Range("D12").Value = -11.7038745880127
$ws.Range("E12").Value = 31.56782150268555
$ws.Range("F12").Value = 1.019708882982564
$ws.Range("G12").Value = -2.124502529607678
$ws.Range("H12").Value = -4.567799980673901

# Row 13
$ws.Range("A13").Value = 1100.0
$ws.Range("B13").Value = "walkingToRunning"
$ws.Range("C13").Value = -1.816783547401428
$ws.Range("D13").Value = 8.313332557678223
$ws.Range("E13").Value = 9.763429641723633
$ws.Range("F13").Value = -3.90443126063479
$ws.Range("G13").Value = -6.289577069683729
$ws.Range("H13").Value = 0.08947662911672927

# Row 14
$ws.Range("A14").Value = 1200.0
$ws.Range("B14").Value = "walkingToRunning"
$ws.Range("C14").Value = 2.05627965927124
$ws.Range("D14").Value = -54.87781143188477
$ws.Range("E14").Value = 29.65373039245605
$ws.Range("F14").Value = -4.737582387211194
$ws.Range("G14").Value = 1.679735990328154
$ws.Range("H14").Value = 5.945494293052519

# Row 15
$ws.Range("A15").Value = 1300.0
$ws.Range("B15").Value = "walkingToRunning"
$ws.Range("C15").Value = -26.46129989624023
$ws.Range("D15").Value = 36.53782653808594
$ws.Range("E15").Value = 2.705925941467285
$ws.Range("F15").Value = -2.958322148456745
$ws.Range("G15").Value = 5.940943116339433
$ws.Range("H15").Value = 1.412405640165503

# Row 16
$ws.Range("A16").Value = 1400.0
$ws.Range("B16").Value = "walkingToRunning"
$ws.Range("C16").Value = -6.459963798522949
$ws.Range("D16").Value = 9.631237030029297
$ws.Range("E16").Value = -4.527087211608887
$ws.Range("F16").Value = -3.080484191948098
$ws.Range("G16").Value = 14.34931452920528
$ws.Range("H16").Value = -6.689891151178712

# Row 17
$ws.Range("A17").Value = 1500.0
$ws.Range("B17").Value = "walkingToRunning"
$ws.Range("C17").Value = -2.323664665222168
$ws.Range("D17").Value = -6.97331714630127
$ws.Range("E17").Value = 3.054933786392212
$ws.Range("F17").Value = 3.206049725274078
$ws.Range("G17").Value = -2.906924434911418
$ws.Range("H17").Value = -3.371373475154916

# Row 18
$ws.Range("A18").Value = 1600.0
$ws.Range("B18").Value = "walkingToRunning"
$ws.Range("C18").Value = 23.92743682861328
$ws.Range("D18").Value = 6.633898258209229
$ws.Range("E18").Value = 20.45568084716797
$ws.Range("F18").Value = 1.706316604792587
$ws.Range("G18").Value = -3.271737626170626
$ws.Range("H18").Value = -8.0454896521347

# Row 19
$ws.Range("A19").Value = 1700.0
$ws.Range("B19").Value = "walkingToRunning"
$ws.Range("C19").Value = 12.88255214691162
$ws.Range("D19").Value = 13.76539325714111
$ws.Range("E19").Value = 5.36094856262207
$ws.Range("F19").Value = -3.738767730855505
$ws.Range("G19").Value = 0.06012086237769707
$ws.Range("H19").Value = 3.632486434740414

# Row 20
$ws.Range("A20").Value = 1800.0
$ws.Range("B20").Value = "walkingToRunning"
$ws.Range("C20").Value = -39.15726470947266
$ws.Range("D20").Value = -50.35159301757812
$ws.Range("E20").Value = 59.05854797363281
$ws.Range("F20").Value = -3.090678747569279
$ws.Range("G20").Value = 1.963554282054647
$ws.Range("H20").Value = 5.732710114149367

# Row 21
$ws.Range("A21").Value = 1900.0
$ws.Range("B21").Value = "walkingToRunning"
$ws.Range("C21").Value = -29.70075225830078
$ws.Range("D21").Value = 18.98210144042969
$ws.Range("E21").Value = -6.946440696716309
$ws.Range("F21").Value = -3.92490029335021
$ws.Range("G21").Value = 8.097008705139087
$ws.Range("H21").Value = -0.6271078586577614

# Row 22
$ws.Range("A22").Value = 2000.0
$ws.Range("B22").Value = "walkingToRunning"
$ws.Range("C22").Value = -3.437598705291748
$ws.Range("D22").Value = 8.67605209350586
$ws.Range("E22").Value = -6.736623287200928
$ws.Range("F22").Value = -1.013831214370013
$ws.Range("G22").Value = 11.75312601071651
$ws.Range("H22").Value = -9.512279387946442

# Row 23
$ws.Range("A23").Value = 2100.0
$ws.Range("B23").Value = "walkingToRunning"
$ws.Range("C23").Value = 36.55035400390625
$ws.Range("D23").Value = -4.483262062072754
$ws.Range("E23").Value = -3.189533472061157
$ws.Range("F23").Value = 1.26375397343502
$ws.Range("G23").Value = -1.259433144720926
$ws.Range("H23").Value = -2.911592367653529

# Row 24
$ws.Range("A24").Value = 2200.0
$ws.Range("B24").Value = "walkingToRunning"
$ws.Range("C24").Value = 25.18490791320801
$ws.Range("D24").Value = 10.66421031951904
$ws.Range("E24").Value = 36.16248321533203
$ws.Range("F24").Value = 3.293022926722733
$ws.Range("G24").Value = -6.641965772504008
$ws.Range("H24").Value = -9.101227296847053

# Row 25
$ws.Range("A25").Value = 2300.0
$ws.Range("B25").Value = "walkingToRunning"
$ws.Range("C25").Value = 11.78367233276367
$ws.Range("D25").Value = 19.32002067565918
$ws.Range("E25").Value = 14.97337532043457
$ws.Range("F25").Value = -0.6729963912028092
$ws.Range("G25").Value = -6.046104440065187
$ws.Range("H25").Value = 2.590652748803134

# Row 26
$ws.Range("A26").Value = 2400.0
$ws.Range("B26").Value = "walkingToRunning"
$ws.Range("C26").Value = -10.2064151763916
$ws.Range("D26").Value = -54.4849967956543
$ws.Range("E26").Value = 45.01205825805664
$ws.Range("F26").Value = -4.304234186622589
$ws.Range("G26").Value = 1.124698318053637
$ws.Range("H26").Value = 5.37724533482133

# Row 27
$ws.Range("A27").Value = 2500.0
$ws.Range("B27").Value = "walkingToRunning"
$ws.Range("C27").Value = -5.304520606994629
$ws.Range("D27").Value = 4.910325050354004
$ws.Range("E27").Value = -39.37523651123047
$ws.Range("F27").Value = -3.525356531143185
$ws.Range("G27").Value = 1.932811015120178
$ws.Range("H27").Value = -0.6847266883493761

# Row 28
$ws.Range("A28").Value = 2600.0
$ws.Range("B28").Value = "walkingToRunning"
$ws.Range("C28").Value = 10.93332672119141
$ws.Range("D28").Value = 14.53017807006836
$ws.Range("E28").Value = -3.883467674255371
$ws.Range("F28").Value = 0.02450294917989604
$ws.Range("G28").Value = 11.26465672644501
$ws.Range("H28").Value = -6.819475414596996

# Row 29
$ws.Range("A29").Value = 2700.0
$ws.Range("B29").Value = "walkingToRunning"
$ws.Range("C29").Value = 23.18131637573243
$ws.Range("D29").Value = -33.343994140625
$ws.Range("E29").Value = -4.427485942840576
$ws.Range("F29").Value = 1.655811092563884
$ws.Range("G29").Value = -3.081682441390621
$ws.Range("H29").Value = -0.3115755630151993

# Row 30
$ws.Range("A30").Value = 2800.0
$ws.Range("B30").Value = "walkingToRunning"
$ws.Range("C30").Value = 17.88149261474609
$ws.Range("D30").Value = -17.56607437133789
$ws.Range("E30").Value = -2.040470600128174
$ws.Range("F30").Value = 4.203615091671459
$ws.Range("G30").Value = -2.449124443196801
$ws.Range("H30").Value = -8.744278280116706

# Row 31
$ws.Range("A31").Value = 2900.0
$ws.Range("B31").Value = "walkingToRunning"
$ws.Range("C31").Value = -4.381585597991943
$ws.Range("D31").Value = 13.65173721313477
$ws.Range("E31").Value = 6.552346229553223
$ws.Range("F31").Value = -2.318640726749061
$ws.Range("G31").Value = -4.338136873512635
$ws.Range("H31").Value = 1.472857043007876

